# Update "Master Inventory" (Sheet1) stock counts and refresh the saved
# view/selection state to match the author's latest working position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Decrement stock quantities in column B for several rows (inventory used
# up / removed crappy FW experimentation parts).
$ws.Range("B9").Value  = 3
$ws.Range("B10").Value = 8
$ws.Range("B32").Value = 10
$ws.Range("B33").Value = 4
$ws.Range("B38").Value = 8
$ws.Range("B40").Value = 13
$ws.Range("B41").Value = 12

# Update the sheet's scroll position and active cell selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select()
